$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '274.43'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.49%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.70'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.52%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.862'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.54%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06326'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.11%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.888'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.38%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.323'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.82%'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.276'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '35.14%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8697'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.06%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1458'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.75%'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.05048'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.63%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07374'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.54%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02923'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-7.28%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09046'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.03%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001573'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.60%'
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006315'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.82%'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.005979'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.01%'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.447'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.13%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.297'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.47%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3153'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.12%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1324'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.98%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.897'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.35%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04361'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.03%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001177'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.21%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004266'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.24%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.07%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001692'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.12%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04040'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.39%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006666'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.85%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1168'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.22%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.77%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01221'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-12.14%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005311'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.36%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.452'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-38.40%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01998'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-33.10%'
